$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5600
$ws.Range("I76").Value = 5600
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5600
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -5285
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 5600
$ws.Range("I79").Value = 5600
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5600
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -4508
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 731.88464
$ws.Range("I98").Value = 681.16
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 681.16
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 816.84
$ws.Range("N98").Value = -4996
$ws.Range("H113").Value = 92126.766
$ws.Range("I113").Value = 213140
$ws.Range("J113").Value = 16493.5
$ws.Range("K113").Value = 213140
$ws.Range("L113").Value = 16493.5
$ws.Range("M113").Value = -209886
$ws.Range("N113").Value = -23001.5
$ws.Range("H122").Value = 731.88464
$ws.Range("I122").Value = 681.16
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2043.48
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 406.52
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 7599536
$ws.Range("I132").Value = 7937540
$ws.Range("K132").Value = 23812620
$ws.Range("M132").Value = -23810090
$ws.Range("H135").Value = 890
$ws.Range("I135").Value = 419.6316
$ws.Range("K135").Value = 3776.6844
$ws.Range("M135").Value = -1241.6844

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 9
$ws.Range("J17").Value = 9
$ws.Range("L17").Value = 9
$ws.Range("N17").Value = -355
$ws.Range("H29").Value = 2470
$ws.Range("I29").Value = 1900
$ws.Range("J29").Value = 2755
$ws.Range("K29").Value = 1900
$ws.Range("L29").Value = 2755
$ws.Range("M29").Value = -1592
$ws.Range("N29").Value = -3371
$ws.Range("H45").Value = 7476.6
$ws.Range("J45").Value = 4216.75
$ws.Range("L45").Value = 4216.75
$ws.Range("N45").Value = -4970.75
$ws.Range("H61").Value = 3421.7144
$ws.Range("I61").Value = 2948.5518
$ws.Range("K61").Value = 2948.5518
$ws.Range("M61").Value = -2736.5518
$ws.Range("H74").Value = 6945.636
$ws.Range("J74").Value = 37539
$ws.Range("L74").Value = 37539
$ws.Range("N74").Value = -39287
$ws.Range("H77").Value = 6945.636
$ws.Range("J77").Value = 37539
$ws.Range("L77").Value = 187695
$ws.Range("N77").Value = -196431
$ws.Range("H136").Value = 3421.7144
$ws.Range("I136").Value = 2948.5518
$ws.Range("K136").Value = 8845.6554
$ws.Range("M136").Value = -6295.6554

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4338.737
$ws.Range("I20").Value = 3545
$ws.Range("J20").Value = 5053.1
$ws.Range("K20").Value = 3545
$ws.Range("L20").Value = 5053.1
$ws.Range("M20").Value = -3298
$ws.Range("N20").Value = -5547.1
$ws.Range("H86").Value = 2540.2307
$ws.Range("I86").Value = 1644.1578
$ws.Range("J86").Value = 4972.4287
$ws.Range("K86").Value = 1644.1578
$ws.Range("L86").Value = 4972.4287
$ws.Range("M86").Value = -521.1578
$ws.Range("N86").Value = -7218.4287
$ws.Range("H89").Value = 2540.2307
$ws.Range("I89").Value = 1644.1578
$ws.Range("J89").Value = 4972.4287
$ws.Range("K89").Value = 8220.789000000001
$ws.Range("L89").Value = 24862.1435
$ws.Range("M89").Value = -2604.789000000001
$ws.Range("N89").Value = -36094.14350000001
$ws.Range("H134").Value = 1903.826
$ws.Range("I134").Value = 1605.7949
$ws.Range("K134").Value = 4817.384700000001
$ws.Range("M134").Value = -2282.384700000001
$ws.Range("H139").Value = 94937.5
$ws.Range("J139").Value = 99916.664
$ws.Range("L139").Value = 99916.664
$ws.Range("N139").Value = -110196.664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6997.8
$ws.Range("J4").Value = 4996.3335
$ws.Range("L4").Value = 4996.3335
$ws.Range("N4").Value = -5220.3335
$ws.Range("H22").Value = 588
$ws.Range("I22").Value = 381.66666
$ws.Range("K22").Value = 381.66666
$ws.Range("M22").Value = -31.66665999999998
$ws.Range("H31").Value = 98517
$ws.Range("J31").Value = 11678.167
$ws.Range("L31").Value = 11678.167
$ws.Range("N31").Value = -12268.167
$ws.Range("H34").Value = 98517
$ws.Range("J34").Value = 11678.167
$ws.Range("L34").Value = 11678.167
$ws.Range("N34").Value = -12082.167
$ws.Range("H105").Value = 1177
$ws.Range("I105").Value = 1333
$ws.Range("K105").Value = 1333
$ws.Range("M105").Value = 414
$ws.Range("H132").Value = 2296.9592
$ws.Range("I132").Value = 2212.2888
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 6636.866399999999
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -4106.866399999999
$ws.Range("N132").Value = -14808.5
$ws.Range("H134").Value = 27796.445
$ws.Range("I134").Value = 17536.834
$ws.Range("K134").Value = 52610.50199999999
$ws.Range("M134").Value = -50075.50199999999
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4674.4
$ws.Range("J116").Value = 4468.1665
$ws.Range("L116").Value = 13404.4995
$ws.Range("N116").Value = -20288.4995
$ws.Range("H131").Value = 42323.48
$ws.Range("I131").Value = 111988.445
$ws.Range("J131").Value = 3136.9375
$ws.Range("K131").Value = 335965.335
$ws.Range("L131").Value = 9410.8125
$ws.Range("M131").Value = -330925.335
$ws.Range("N131").Value = -19490.8125
$ws.Range("H139").Value = 3699.7
$ws.Range("I139").Value = 2025
$ws.Range("K139").Value = 6075
$ws.Range("M139").Value = -935

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10033.333
$ws.Range("I5").Value = 5050
$ws.Range("J5").Value = 20000
$ws.Range("K5").Value = 5050
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = -4938
$ws.Range("N5").Value = -20224
$ws.Range("H70").Value = 12180.733
$ws.Range("J70").Value = 14106.777
$ws.Range("L70").Value = 14106.777
$ws.Range("N70").Value = -14646.777
$ws.Range("H73").Value = 12180.733
$ws.Range("J73").Value = 14106.777
$ws.Range("L73").Value = 14106.777
$ws.Range("N73").Value = -15978.777
$ws.Range("H132").Value = 3192.3928
$ws.Range("I132").Value = 3027.7
$ws.Range("K132").Value = 9083.099999999999
$ws.Range("M132").Value = -6553.099999999999
$ws.Range("H135").Value = 80780
$ws.Range("J135").Value = 80780
$ws.Range("L135").Value = 80780
$ws.Range("N135").Value = -90920

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 11368848
$ws.Range("I16").Value = 26318040
$ws.Range("J16").Value = 7462.2
$ws.Range("K16").Value = 26318040
$ws.Range("L16").Value = 7462.2
$ws.Range("M16").Value = -26317870
$ws.Range("N16").Value = -7802.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 42856.57
$ws.Range("J2").Value = 49999.25
$ws.Range("L2").Value = 49999.25
$ws.Range("N2").Value = -50223.25
$ws.Range("H12").Value = 2988
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H81").Value = 8849.65
$ws.Range("I81").Value = 16015.143
$ws.Range("J81").Value = 4991.3076
$ws.Range("K81").Value = 32030.286
$ws.Range("L81").Value = 9982.6152
$ws.Range("M81").Value = -30969.286
$ws.Range("N81").Value = -12104.6152
$ws.Range("H84").Value = 8849.65
$ws.Range("I84").Value = 16015.143
$ws.Range("J84").Value = 4991.3076
$ws.Range("K84").Value = 160151.43
$ws.Range("L84").Value = 49913.076
$ws.Range("M84").Value = -154847.43
$ws.Range("N84").Value = -60521.076
$ws.Range("H136").Value = 1690.5667
$ws.Range("I136").Value = 1673.8334
$ws.Range("K136").Value = 5021.5002
$ws.Range("M136").Value = -2471.5002
